$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new row 44 with the latest Kaspa buy entry (run on 2025-07-25).
# Force column A to be treated as plain text so the date string
# "07/25/2025" is preserved instead of being auto-converted to a date
# serial number, then copy the (unstyled) format of the preceding data
# row so no stray style index is left on the new cell.
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = "07/25/2025"
$ws.Cells.Item(44, 1).Style = $ws.Cells.Item(43, 1).Style

$ws.Cells.Item(44, 2).Value = 499.4599999999991
$ws.Cells.Item(44, 3).Value = 0.1001081167661076
$ws.Cells.Item(44, 4).Value = 50
